$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# The commit removes the two data rows that previously sat at rows 2 and 3
# (Hydnellum gracilipes / Anomoporia kamtschatica records), shifting every
# following record up by two rows. Deleting the entire rows reproduces the
# observed result: the sheet's used range shrinks from A1:AY10 to A1:AY8 and
# every remaining record (previously rows 4-10) now occupies rows 2-8 with
# all of its original cell values intact.
$ws.Range("A2:A3").EntireRow.Delete()
